$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet: "Property1" -> "DataNode" (unify DataNode/DataTable/Entity naming)
$ws.Name = "DataNode"

# Header row (row 1) now wraps to two lines -> taller row height
$ws.Rows.Item(1).RowHeight = 27

# Row 8 (the "Desc" row) shrinks slightly
$ws.Rows.Item(8).RowHeight = 81

# Leave the last active selection on H13, matching the saved UI state
$ws.Range("H13").Select()
